$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.586.60"
$ws.Cells.Item(2, 5).Value = "  +2.21%  "
$ws.Cells.Item(3, 4).Value = "2.387.77"
$ws.Cells.Item(3, 5).Value = "  +2.19%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'553.23"
$ws.Cells.Item(5, 5).Value = "  +2.13%  "
$ws.Cells.Item(6, 4).Value = "'140.61"
$ws.Cells.Item(6, 5).Value = "  +2.60%  "
$ws.Cells.Item(7, 5).Value = "  -0.06%  "
$ws.Cells.Item(8, 4).Value = "'0.527"
$ws.Cells.Item(8, 5).Value = "  +1.20%  "
$ws.Cells.Item(9, 4).Value = "2.389.46"
$ws.Cells.Item(9, 5).Value = "  +2.33%  "
$ws.Cells.Item(10, 4).Value = "'0.109"
$ws.Cells.Item(10, 5).Value = "  +5.52%  "
$ws.Cells.Item(11, 5).Value = "  +1.85%  "
$ws.Cells.Item(12, 4).Value = "'5.36"
$ws.Cells.Item(12, 5).Value = "  +3.52%  "
$ws.Cells.Item(13, 4).Value = "'0.352"
$ws.Cells.Item(13, 5).Value = "  +4.72%  "
$ws.Cells.Item(14, 4).Value = "'25.53"
$ws.Cells.Item(14, 5).Value = "  +3.36%  "
$ws.Cells.Item(15, 4).Value = "'0.0000168"
$ws.Cells.Item(15, 5).Value = "  +4.88%  "
$ws.Cells.Item(16, 4).Value = "61.499.98"
$ws.Cells.Item(16, 5).Value = "  +1.22%  "
$ws.Cells.Item(17, 2).Value = "Chainlink"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(17, 4).Value = "'11.04"
$ws.Cells.Item(17, 5).Value = "  +5.24%  "
$ws.Cells.Item(18, 2).Value = "BitcoinCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(18, 4).Value = "'322.34"
$ws.Cells.Item(18, 5).Value = "  +2.99%  "
$ws.Cells.Item(19, 2).Value = "Polkadot"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(19, 4).Value = "'4.16"
$ws.Cells.Item(19, 5).Value = "  +2.60%  "
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Value = "'6.80"
$ws.Cells.Item(20, 5).Value = "  +4.87%  "
$ws.Cells.Item(21, 2).Value = "Dai"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(21, 4).Value = "'1.00"
$ws.Cells.Item(21, 5).Value = "  +0.19%  "
$ws.Cells.Item(22, 2).Value = "Litecoin"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(22, 4).Value = "'64.50"
$ws.Cells.Item(22, 5).Value = "  +2.83%  "
$ws.Cells.Item(23, 2).Value = "SuiNetwork"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(23, 4).Value = "'1.74"
$ws.Cells.Item(23, 5).Value = "  -5.86%  "
$ws.Cells.Item(24, 2).Value = "Aptos"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(24, 4).Value = "'8.84"
$ws.Cells.Item(24, 5).Value = "  +10.41%  "
$ws.Cells.Item(25, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(25, 4).Value = "'8.28"
$ws.Cells.Item(25, 5).Value = "  +4.85%  "
$ws.Cells.Item(26, 2).Value = "Bittensor"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(26, 4).Value = "'525.17"
$ws.Cells.Item(26, 5).Value = "  +4.01%  "
$ws.Cells.Item(27, 2).Value = "PEPE"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(27, 4).Value = "0.0₃0911"
$ws.Cells.Item(27, 5).Value = "  +2.10%  "
$ws.Cells.Item(28, 2).Value = "Kaspa"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(28, 4).Value = "'0.151"
$ws.Cells.Item(28, 5).Value = "  +5.53%  "
$ws.Cells.Item(29, 2).Value = "Fetch.AI"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(29, 4).Value = "'1.39"
$ws.Cells.Item(29, 5).Value = "  +1.42%  "
$ws.Cells.Item(30, 2).Value = "PancakeSwap"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30, 4).Value = "'1.85"
$ws.Cells.Item(30, 5).Value = "  +3.17%  "
$ws.Cells.Item(31, 2).Value = "ImmutableX"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(31, 4).Value = "'1.56"
$ws.Cells.Item(31, 5).Value = "  +2.37%  "
$ws.Cells.Item(32, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(32, 4).Value = "'0.999"
$ws.Cells.Item(32, 5).Value = "  +0.04%  "
$ws.Cells.Item(33, 2).Value = "NEARProtocol"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(33, 4).Value = "'4.74"
$ws.Cells.Item(33, 5).Value = "  +5.12%  "
$ws.Cells.Item(34, 2).Value = "RenderToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(34, 4).Value = "'5.55"
$ws.Cells.Item(34, 5).Value = "  +6.69%  "
$ws.Cells.Item(35, 2).Value = "Stacks"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(35, 4).Value = "'1.92"
$ws.Cells.Item(35, 5).Value = "  +8.44%  "
$ws.Cells.Item(36, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(36, 4).Value = "'0.380"
$ws.Cells.Item(36, 5).Value = "  +2.81%  "
$ws.Cells.Item(37, 2).Value = "EthereumClassic"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(37, 4).Value = "'18.60"
$ws.Cells.Item(37, 5).Value = "  +2.67%  "
$ws.Cells.Item(38, 2).Value = "Monero"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(38, 4).Value = "'147.29"
$ws.Cells.Item(38, 5).Value = "  +6.11%  "
$ws.Cells.Item(39, 2).Value = "USDe"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(39, 4).Value = "'0.999"
$ws.Cells.Item(39, 5).Value = "  -0.11%  "
$ws.Cells.Item(40, 2).Value = "OKB"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(40, 4).Value = "'41.41"
$ws.Cells.Item(40, 5).Value = "  +3.64%  "
$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 4).Value = "'150.17"
$ws.Cells.Item(41, 5).Value = "  +10.09%  "
$ws.Cells.Item(42, 2).Value = "dogwifhat"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(42, 4).Value = "'2.18"
$ws.Cells.Item(42, 5).Value = "  +5.29%  "
$ws.Cells.Item(43, 2).Value = "Filecoin"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(43, 4).Value = "'3.62"
$ws.Cells.Item(43, 5).Value = "  +2.94%  "
$ws.Cells.Item(44, 2).Value = "Hedera"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(44, 4).Value = "'0.0526"
$ws.Cells.Item(44, 5).Value = "  +3.29%  "
$ws.Cells.Item(45, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(45, 4).Value = "'19.82"
$ws.Cells.Item(45, 5).Value = "  +2.45%  "
$ws.Cells.Item(46, 2).Value = "Mantle"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(46, 4).Value = "'0.584"
$ws.Cells.Item(46, 5).Value = "  +3.78%  "
$ws.Cells.Item(47, 2).Value = "Stellar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(47, 4).Value = "'0.0910"
$ws.Cells.Item(47, 5).Value = "  +1.85%  "
$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).Value = "'0.0225"
$ws.Cells.Item(48, 5).Value = "  +2.18%  "
$ws.Cells.Item(49, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(49, 4).Value = "'11.42"
$ws.Cells.Item(49, 5).Value = "  +0.85%  "
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "'16.87"
$ws.Cells.Item(50, 5).Value = "  +2.67%  "
$ws.Cells.Item(51, 2).Value = "BitgetToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Cells.Item(51, 4).Value = "'1.01"
$ws.Cells.Item(51, 5).Value = "  +4.50%  "
